$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 3: add PriceChange (X3) and UpDown (Y3) values
$ws.Range("X3").Value = -0.29999499999999557
$ws.Range("Y3").Value = "Down"

# Add new row 4 of data
$ws.Range("A4").Value = 42635.817361111112
$ws.Range("A4").NumberFormat = "m/d/yy h:mm"
$ws.Range("B4").Value = -4
$ws.Range("C4").Value = "Neutral"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = "Random"
$ws.Range("Q4").Value = 63.486785924529997
$ws.Range("R4").Value = 1.76
$ws.Range("S4").Value = 0.109
$ws.Range("S4").NumberFormat = "0.00%"
$ws.Range("T4").Value = 0.0455
$ws.Range("T4").NumberFormat = "0.00%"
$ws.Range("U4").Value = 4.84
$ws.Range("V4").Value = 2.2799999999999998
$ws.Range("W4").Value = 0

# Column C width adjustment (widened slightly to fit "Neutral")
$ws.Columns("C:C").ColumnWidth = 6.8
